# Auto-generated Excel COM-interop script applying numeric cell updates
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 21218.75
$ws.Range("I51").Value = 11083.333
$ws.Range("J51").Value = 27300
$ws.Range("K51").Value = 11083.333
$ws.Range("L51").Value = 27300
$ws.Range("M51").Value = -10599.333
$ws.Range("N51").Value = -28268
$ws.Range("H88").Value = 575.9091
$ws.Range("J88").Value = 581.1111
$ws.Range("L88").Value = 581.1111
$ws.Range("N88").Value = -1393.1111
$ws.Range("H91").Value = 575.9091
$ws.Range("J91").Value = 581.1111
$ws.Range("L91").Value = 581.1111
$ws.Range("N91").Value = -3389.1111
$ws.Range("H99").Value = 338.7143
$ws.Range("I99").Value = 316.2
$ws.Range("J99").Value = 395
$ws.Range("K99").Value = 948.5999999999999
$ws.Range("L99").Value = 1185
$ws.Range("M99").Value = 549.4000000000001
$ws.Range("N99").Value = -4181
$ws.Range("H137").Value = 9619850
$ws.Range("I137").Value = 25006404
$ws.Range("K137").Value = 75019212
$ws.Range("M137").Value = -75016662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7080123
$ws.Range("I2").Value = 15336050
$ws.Range("J2").Value = 3614.7144
$ws.Range("K2").Value = 15336050
$ws.Range("L2").Value = 3614.7144
$ws.Range("M2").Value = -15335937
$ws.Range("N2").Value = -3840.7144
$ws.Range("H8").Value = 4316.857
$ws.Range("I8").Value = 5600
$ws.Range("K8").Value = 5600
$ws.Range("M8").Value = -5456
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H61").Value = 6225.6113
$ws.Range("I61").Value = 3393.7144
$ws.Range("K61").Value = 3393.7144
$ws.Range("M61").Value = -3181.7144
$ws.Range("H94").Value = 45972
$ws.Range("J94").Value = 45972
$ws.Range("L94").Value = 45972
$ws.Range("N94").Value = -47774
$ws.Range("H97").Value = 2533091
$ws.Range("I97").Value = 2649739.5
$ws.Range("J97").Value = 900011
$ws.Range("K97").Value = 2649739.5
$ws.Range("L97").Value = 900011
$ws.Range("M97").Value = -2649243.5
$ws.Range("N97").Value = -901003
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H110").Value = 10000914
$ws.Range("I110").Value = 11364312
$ws.Range("K110").Value = 11364312
$ws.Range("M110").Value = -11362267
$ws.Range("H116").Value = 7080123
$ws.Range("I116").Value = 15336050
$ws.Range("J116").Value = 3614.7144
$ws.Range("K116").Value = 15336050
$ws.Range("L116").Value = 3614.7144
$ws.Range("M116").Value = -15333756
$ws.Range("N116").Value = -8202.7144
$ws.Range("H122").Value = 1832.1666
$ws.Range("I122").Value = 1832.1666
$ws.Range("K122").Value = 5496.4998
$ws.Range("M122").Value = -3046.4998
$ws.Range("H132").Value = 3810.6428
$ws.Range("I132").Value = 3328.32
$ws.Range("K132").Value = 9984.960000000001
$ws.Range("M132").Value = -7454.960000000001
$ws.Range("H136").Value = 6225.6113
$ws.Range("I136").Value = 3393.7144
$ws.Range("K136").Value = 10181.1432
$ws.Range("M136").Value = -7631.143199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7080123
$ws.Range("I3").Value = 15336050
$ws.Range("J3").Value = 3614.7144
$ws.Range("K3").Value = 15336050
$ws.Range("L3").Value = 3614.7144
$ws.Range("M3").Value = -15335936
$ws.Range("N3").Value = -3842.7144
$ws.Range("H80").Value = 80301.305
$ws.Range("J80").Value = 104192.4
$ws.Range("L80").Value = 104192.4
$ws.Range("N80").Value = -106188.4
$ws.Range("H83").Value = 80301.305
$ws.Range("J83").Value = 104192.4
$ws.Range("L83").Value = 520962
$ws.Range("N83").Value = -530946
$ws.Range("H86").Value = 184098.9
$ws.Range("I86").Value = 2589.3333
$ws.Range("J86").Value = 1000892
$ws.Range("K86").Value = 2589.3333
$ws.Range("L86").Value = 1000892
$ws.Range("M86").Value = -1466.3333
$ws.Range("N86").Value = -1003138
$ws.Range("H89").Value = 184098.9
$ws.Range("I89").Value = 2589.3333
$ws.Range("J89").Value = 1000892
$ws.Range("K89").Value = 12946.6665
$ws.Range("L89").Value = 5004460
$ws.Range("M89").Value = -7330.666499999999
$ws.Range("N89").Value = -5015692
$ws.Range("H103").Value = 54219
$ws.Range("J103").Value = 54219
$ws.Range("L103").Value = 54219
$ws.Range("N103").Value = -56563
$ws.Range("H112").Value = 98465.664
$ws.Range("J112").Value = 98465.664
$ws.Range("L112").Value = 98465.664
$ws.Range("N112").Value = -101419.664
$ws.Range("H115").Value = 99995
$ws.Range("J115").Value = 99995
$ws.Range("L115").Value = 99995
$ws.Range("N115").Value = -103129
$ws.Range("H134").Value = 3786.7046
$ws.Range("I134").Value = 1330.0714
$ws.Range("K134").Value = 3990.2142
$ws.Range("M134").Value = -1455.2142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2371.077
$ws.Range("I105").Value = 1894.25
$ws.Range("J105").Value = 3134
$ws.Range("K105").Value = 1894.25
$ws.Range("L105").Value = 3134
$ws.Range("M105").Value = -147.25
$ws.Range("N105").Value = -6628
$ws.Range("H107").Value = 451.2
$ws.Range("I107").Value = 439
$ws.Range("K107").Value = 439
$ws.Range("M107").Value = 1481
$ws.Range("H134").Value = 4764.85
$ws.Range("I134").Value = 4504.6763
$ws.Range("J134").Value = 6239.1665
$ws.Range("K134").Value = 13514.0289
$ws.Range("L134").Value = 18717.4995
$ws.Range("M134").Value = -10979.0289
$ws.Range("N134").Value = -23787.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 89048.22
$ws.Range("J37").Value = 89048.22
$ws.Range("L37").Value = 267144.66
$ws.Range("N37").Value = -267368.66
$ws.Range("H55").Value = 1972.9231
$ws.Range("I55").Value = 709.8
$ws.Range("J55").Value = 2762.375
$ws.Range("K55").Value = 2129.4
$ws.Range("L55").Value = 8287.125
$ws.Range("M55").Value = -1952.4
$ws.Range("N55").Value = -8641.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H97").Value = 1330.9231
$ws.Range("I97").Value = 1381.4166
$ws.Range("K97").Value = 1381.4166
$ws.Range("M97").Value = -885.4166
$ws.Range("H107").Value = 756
$ws.Range("I107").Value = 807.2
$ws.Range("K107").Value = 807.2
$ws.Range("M107").Value = 1112.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 12000
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H55").Value = 24306.477
$ws.Range("I55").Value = 50670.5
$ws.Range("J55").Value = 339.18182
$ws.Range("K55").Value = 50670.5
$ws.Range("L55").Value = 339.18182
$ws.Range("M55").Value = -50497.5
$ws.Range("N55").Value = -685.18182
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H93").Value = 1882
$ws.Range("I93").Value = 1457.1428
$ws.Range("J93").Value = 2625.5
$ws.Range("K93").Value = 1457.1428
$ws.Range("L93").Value = 2625.5
$ws.Range("M93").Value = -209.1428000000001
$ws.Range("N93").Value = -5121.5
$ws.Range("H136").Value = 5315.0435
$ws.Range("J136").Value = 11308.167
$ws.Range("L136").Value = 33924.501
$ws.Range("N136").Value = -39024.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 68974
$ws.Range("J131").Value = 68974
$ws.Range("L131").Value = 68974
$ws.Range("N131").Value = -79054
